$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell so the dataset is R-friendly (snake_case column name)
$ws.Range("A1").Value = "sale_amount"

# Drop the currency formatting on the data column in favour of a plain
# numeric format (no "$" symbol) - also R-friendly for read.csv/readxl etc.
$ws.Columns("A:A").NumberFormat = "0.00"

# Select the whole column (matches the saved selection state)
$ws.Columns("A:A").Select() | Out-Null
